$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay as text (avoid numeric auto-conversion/locale parsing)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.395.36"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "1.937.88"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "0.7706"
$ws.Range("E5").Value = "  +6.36%  "
$ws.Range("D6").Value = "248.82"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "28.14"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "0.3208"
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").Value = "0.07108"
$ws.Range("E10").Value = "  -2.52%  "
$ws.Range("D11").Value = "0.7886"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "0.07996"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "1.936.79"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "5.386"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "94.86"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  -3.32%  "
$ws.Range("D17").Value = "30.397.75"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "256.38"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "0.000008034"
$ws.Range("E19").Value = "  -3.17%  "
$ws.Range("D20").Value = "5.807"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "2.189.45"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "6.813"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "9.617"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "164.78"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "0.1353"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").Value = "19.17"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").Value = "2.301"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").Value = "1.376"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").Value = "4.445"
$ws.Range("D33").Value = "4.162"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "0.05208"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "1.287"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").Value = "0.7520"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "2.777"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "0.01978"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "2.810"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "78.41"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "6.477"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("D42").Value = "0.4538"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "1.989"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "101.93"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").Value = "7.569"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").Value = "9.834"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "37.68"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "979.98"
$ws.Range("E50").Value = "  +10.82%  "
$ws.Range("D51").Value = "0.4178"
$ws.Range("E51").Value = "  -0.36%  "
